$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 42   # H11: 39 -> 42
$ws.Cells.Item(11, 9).Value = 42   # I11: 39 -> 42
$ws.Cells.Item(11, 11).Value = 42   # K11: 39 -> 42
$ws.Cells.Item(11, 13).Value = 98   # M11: 101 -> 98
$ws.Cells.Item(40, 8).Value = 26670.223   # H40: 29874 -> 26670.223
$ws.Cells.Item(40, 9).Value = 27009.25   # I40: 35665.668 -> 27009.25
$ws.Cells.Item(40, 11).Value = 27009.25   # K40: 35665.668 -> 27009.25
$ws.Cells.Item(40, 13).Value = -26834.25   # M40: -35490.668 -> -26834.25
$ws.Cells.Item(93, 8).Value = 20000   # H93: 0 -> 20000
$ws.Cells.Item(93, 10).Value = 20000   # J93: 0 -> 20000
$ws.Cells.Item(93, 12).Value = 20000   # L93: 0 -> 20000
$ws.Cells.Item(93, 14).Value = -24992   # N93: None -> -24992
$ws.Cells.Item(100, 8).Value = 2010.7646   # H100: 2393.7856 -> 2010.7646
$ws.Cells.Item(100, 9).Value = 1600.375   # I100: 1797.5714 -> 1600.375
$ws.Cells.Item(100, 10).Value = 2375.5557   # J100: 2990 -> 2375.5557
$ws.Cells.Item(100, 11).Value = 1600.375   # K100: 1797.5714 -> 1600.375
$ws.Cells.Item(100, 12).Value = 2375.5557   # L100: 2990 -> 2375.5557
$ws.Cells.Item(100, 13).Value = -1059.375   # M100: -1256.5714 -> -1059.375
$ws.Cells.Item(100, 14).Value = -3457.5557   # N100: -4072 -> -3457.5557
$ws.Cells.Item(116, 8).Value = 31486156   # H116: 24641954 -> 31486156
$ws.Cells.Item(116, 9).Value = 62969136   # I116: 43594890 -> 62969136
$ws.Cells.Item(116, 10).Value = 3174   # J116: 3139.1 -> 3174
$ws.Cells.Item(116, 11).Value = 62969136   # K116: 43594890 -> 62969136
$ws.Cells.Item(116, 12).Value = 3174   # L116: 3139.1 -> 3174
$ws.Cells.Item(116, 13).Value = -62965694   # M116: -43591448 -> -62965694
$ws.Cells.Item(116, 14).Value = -10058   # N116: -10023.1 -> -10058
$ws.Cells.Item(138, 8).Value = 6142.155   # H138: 6105.347 -> 6142.155
$ws.Cells.Item(138, 10).Value = 7152.1523   # J138: 7091.15 -> 7152.1523
$ws.Cells.Item(138, 12).Value = 21456.4569   # L138: 21273.45 -> 21456.4569
$ws.Cells.Item(138, 14).Value = -31736.4569   # N138: -31553.45 -> -31736.4569

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 27293.525   # H2: 28809.834 -> 27293.525
$ws.Cells.Item(2, 9).Value = 32282.125   # I2: 34434.266 -> 32282.125
$ws.Cells.Item(2, 11).Value = 32282.125   # K2: 34434.266 -> 32282.125
$ws.Cells.Item(2, 13).Value = -32169.125   # M2: -34321.266 -> -32169.125
$ws.Cells.Item(32, 8).Value = 16033.533   # H32: 15771.819 -> 16033.533
$ws.Cells.Item(32, 9).Value = 15538.77   # I32: 15246.887 -> 15538.77
$ws.Cells.Item(32, 11).Value = 15538.77   # K32: 15246.887 -> 15538.77
$ws.Cells.Item(32, 13).Value = -15251.77   # M32: -14959.887 -> -15251.77
$ws.Cells.Item(74, 8).Value = 5319826   # H74: 5319839 -> 5319826
$ws.Cells.Item(74, 9).Value = 6250580   # I74: 6250595 -> 6250580
$ws.Cells.Item(74, 11).Value = 6250580   # K74: 6250595 -> 6250580
$ws.Cells.Item(74, 13).Value = -6249706   # M74: -6249721 -> -6249706
$ws.Cells.Item(77, 8).Value = 5319826   # H77: 5319839 -> 5319826
$ws.Cells.Item(77, 9).Value = 6250580   # I77: 6250595 -> 6250580
$ws.Cells.Item(77, 11).Value = 31252900   # K77: 31252975 -> 31252900
$ws.Cells.Item(77, 13).Value = -31248532   # M77: -31248607 -> -31248532
$ws.Cells.Item(96, 8).Value = 27497.5   # H96: 29995 -> 27497.5
$ws.Cells.Item(96, 10).Value = 27497.5   # J96: 29995 -> 27497.5
$ws.Cells.Item(96, 12).Value = 27497.5   # L96: 29995 -> 27497.5
$ws.Cells.Item(96, 14).Value = -32989.5   # N96: -35487 -> -32989.5
$ws.Cells.Item(116, 8).Value = 27293.525   # H116: 28809.834 -> 27293.525
$ws.Cells.Item(116, 9).Value = 32282.125   # I116: 34434.266 -> 32282.125
$ws.Cells.Item(116, 11).Value = 32282.125   # K116: 34434.266 -> 32282.125
$ws.Cells.Item(116, 13).Value = -29988.125   # M116: -32140.266 -> -29988.125

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 27293.525   # H3: 28809.834 -> 27293.525
$ws.Cells.Item(3, 9).Value = 32282.125   # I3: 34434.266 -> 32282.125
$ws.Cells.Item(3, 11).Value = 32282.125   # K3: 34434.266 -> 32282.125
$ws.Cells.Item(3, 13).Value = -32168.125   # M3: -34320.266 -> -32168.125
$ws.Cells.Item(80, 8).Value = 847.5454999999999   # H80: 1020.7 -> 847.5454999999999
$ws.Cells.Item(80, 9).Value = 796.3333   # I80: 799.6667 -> 796.3333
$ws.Cells.Item(80, 10).Value = 866.75   # J80: 1115.4286 -> 866.75
$ws.Cells.Item(80, 11).Value = 796.3333   # K80: 799.6667 -> 796.3333
$ws.Cells.Item(80, 12).Value = 866.75   # L80: 1115.4286 -> 866.75
$ws.Cells.Item(80, 13).Value = 201.6667   # M80: 198.3333 -> 201.6667
$ws.Cells.Item(80, 14).Value = -2862.75   # N80: -3111.4286 -> -2862.75
$ws.Cells.Item(83, 8).Value = 847.5454999999999   # H83: 1020.7 -> 847.5454999999999
$ws.Cells.Item(83, 9).Value = 796.3333   # I83: 799.6667 -> 796.3333
$ws.Cells.Item(83, 10).Value = 866.75   # J83: 1115.4286 -> 866.75
$ws.Cells.Item(83, 11).Value = 3981.6665   # K83: 3998.3335 -> 3981.6665
$ws.Cells.Item(83, 12).Value = 4333.75   # L83: 5577.143 -> 4333.75
$ws.Cells.Item(83, 13).Value = 1010.3335   # M83: 993.6665000000003 -> 1010.3335
$ws.Cells.Item(83, 14).Value = -14317.75   # N83: -15561.143 -> -14317.75

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1911.033   # H31: 1898.3695 -> 1911.033
$ws.Cells.Item(31, 10).Value = 5974.857   # J31: 5321.25 -> 5974.857
$ws.Cells.Item(31, 12).Value = 5974.857   # L31: 5321.25 -> 5974.857
$ws.Cells.Item(31, 14).Value = -6564.857   # N31: -5911.25 -> -6564.857
$ws.Cells.Item(34, 8).Value = 1911.033   # H34: 1898.3695 -> 1911.033
$ws.Cells.Item(34, 10).Value = 5974.857   # J34: 5321.25 -> 5974.857
$ws.Cells.Item(34, 12).Value = 5974.857   # L34: 5321.25 -> 5974.857
$ws.Cells.Item(34, 14).Value = -6378.857   # N34: -5725.25 -> -6378.857
$ws.Cells.Item(141, 8).Value = 118225.02   # H141: 121075.125 -> 118225.02
$ws.Cells.Item(141, 10).Value = 119162.85   # J141: 122207.2 -> 119162.85
$ws.Cells.Item(141, 12).Value = 119162.85   # L141: 122207.2 -> 119162.85
$ws.Cells.Item(141, 14).Value = -129522.85   # N141: -132567.2 -> -129522.85

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(141, 8).Value = 5768.6665   # H141: 5555.5386 -> 5768.6665
$ws.Cells.Item(141, 9).Value = 2817.8572   # I141: 2840.375 -> 2817.8572
$ws.Cells.Item(141, 11).Value = 8453.571599999999   # K141: 8521.125 -> 8453.571599999999
$ws.Cells.Item(141, 13).Value = -3273.571599999999   # M141: -3341.125 -> -3273.571599999999

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7183.1665   # H70: 7673.1665 -> 7183.1665
$ws.Cells.Item(70, 9).Value = 6860.35   # I70: 7263.684 -> 6860.35
$ws.Cells.Item(70, 10).Value = 7828.8   # J70: 8380.454 -> 7828.8
$ws.Cells.Item(70, 11).Value = 6860.35   # K70: 7263.684 -> 6860.35
$ws.Cells.Item(70, 12).Value = 7828.8   # L70: 8380.454 -> 7828.8
$ws.Cells.Item(70, 13).Value = -6590.35   # M70: -6993.684 -> -6590.35
$ws.Cells.Item(70, 14).Value = -8368.799999999999   # N70: -8920.454 -> -8368.799999999999
$ws.Cells.Item(73, 8).Value = 7183.1665   # H73: 7673.1665 -> 7183.1665
$ws.Cells.Item(73, 9).Value = 6860.35   # I73: 7263.684 -> 6860.35
$ws.Cells.Item(73, 10).Value = 7828.8   # J73: 8380.454 -> 7828.8
$ws.Cells.Item(73, 11).Value = 6860.35   # K73: 7263.684 -> 6860.35
$ws.Cells.Item(73, 12).Value = 7828.8   # L73: 8380.454 -> 7828.8
$ws.Cells.Item(73, 13).Value = -5924.35   # M73: -6327.684 -> -5924.35
$ws.Cells.Item(73, 14).Value = -9700.799999999999   # N73: -10252.454 -> -9700.799999999999
$ws.Cells.Item(95, 8).Value = 26445   # H95: 28052.5 -> 26445
$ws.Cells.Item(95, 10).Value = 26445   # J95: 28052.5 -> 26445
$ws.Cells.Item(95, 12).Value = 26445   # L95: 28052.5 -> 26445
$ws.Cells.Item(95, 14).Value = -31937   # N95: -33544.5 -> -31937
$ws.Cells.Item(98, 8).Value = 29624.75   # H98: 31397.8 -> 29624.75
$ws.Cells.Item(98, 10).Value = 29624.75   # J98: 31397.8 -> 29624.75
$ws.Cells.Item(98, 12).Value = 29624.75   # L98: 31397.8 -> 29624.75
$ws.Cells.Item(98, 14).Value = -35614.75   # N98: -37387.8 -> -35614.75
$ws.Cells.Item(132, 8).Value = 6191.8   # H132: 6996.6665 -> 6191.8
$ws.Cells.Item(132, 9).Value = 4986.3335   # I132: 4990 -> 4986.3335
$ws.Cells.Item(132, 11).Value = 14959.0005   # K132: 14970 -> 14959.0005
$ws.Cells.Item(132, 13).Value = -12429.0005   # M132: -12440 -> -12429.0005

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4041.9023   # H7: 4074.1482 -> 4041.9023
$ws.Cells.Item(7, 9).Value = 3432.0178   # I7: 3449.8704 -> 3432.0178
$ws.Cells.Item(7, 10).Value = 5355.5   # J7: 5322.7036 -> 5355.5
$ws.Cells.Item(7, 11).Value = 3432.0178   # K7: 3449.8704 -> 3432.0178
$ws.Cells.Item(7, 12).Value = 5355.5   # L7: 5322.7036 -> 5355.5
$ws.Cells.Item(7, 13).Value = -3320.0178   # M7: -3337.8704 -> -3320.0178
$ws.Cells.Item(7, 14).Value = -5579.5   # N7: -5546.7036 -> -5579.5
$ws.Cells.Item(16, 8).Value = 1525.8182   # H16: 1598.4286 -> 1525.8182
$ws.Cells.Item(16, 9).Value = 1213.0667   # I16: 1162.125 -> 1213.0667
$ws.Cells.Item(16, 10).Value = 2196   # J16: 2994.6 -> 2196
$ws.Cells.Item(16, 11).Value = 1213.0667   # K16: 1162.125 -> 1213.0667
$ws.Cells.Item(16, 12).Value = 2196   # L16: 2994.6 -> 2196
$ws.Cells.Item(16, 13).Value = -1043.0667   # M16: -992.125 -> -1043.0667
$ws.Cells.Item(16, 14).Value = -2536   # N16: -3334.6 -> -2536
$ws.Cells.Item(55, 8).Value = 563.2727   # H55: 527.0833 -> 563.2727
$ws.Cells.Item(55, 9).Value = 103.4   # I55: 107.666664 -> 103.4
$ws.Cells.Item(55, 11).Value = 103.4   # K55: 107.666664 -> 103.4
$ws.Cells.Item(55, 13).Value = 69.59999999999999   # M55: 65.333336 -> 69.59999999999999
$ws.Cells.Item(68, 8).Value = 4032.5588   # H68: 3942.6365 -> 4032.5588
$ws.Cells.Item(68, 10).Value = 7790.857   # J68: 7922.6665 -> 7790.857
$ws.Cells.Item(68, 12).Value = 7790.857   # L68: 7922.6665 -> 7790.857
$ws.Cells.Item(68, 14).Value = -9288.857   # N68: -9420.666499999999 -> -9288.857
$ws.Cells.Item(71, 8).Value = 4032.5588   # H71: 3942.6365 -> 4032.5588
$ws.Cells.Item(71, 10).Value = 7790.857   # J71: 7922.6665 -> 7790.857
$ws.Cells.Item(71, 12).Value = 38954.285   # L71: 39613.3325 -> 38954.285
$ws.Cells.Item(71, 14).Value = -46442.285   # N71: -47101.3325 -> -46442.285
$ws.Cells.Item(82, 8).Value = 1406.7826   # H82: 1556.6957 -> 1406.7826
$ws.Cells.Item(82, 9).Value = 1033.5385   # I82: 1231 -> 1033.5385
$ws.Cells.Item(82, 10).Value = 1892   # J82: 2063.3333 -> 1892
$ws.Cells.Item(82, 11).Value = 1033.5385   # K82: 1231 -> 1033.5385
$ws.Cells.Item(82, 12).Value = 1892   # L82: 2063.3333 -> 1892
$ws.Cells.Item(82, 13).Value = -672.5385000000001   # M82: -870 -> -672.5385000000001
$ws.Cells.Item(82, 14).Value = -2614   # N82: -2785.3333 -> -2614
$ws.Cells.Item(85, 8).Value = 1406.7826   # H85: 1556.6957 -> 1406.7826
$ws.Cells.Item(85, 9).Value = 1033.5385   # I85: 1231 -> 1033.5385
$ws.Cells.Item(85, 10).Value = 1892   # J85: 2063.3333 -> 1892
$ws.Cells.Item(85, 11).Value = 1033.5385   # K85: 1231 -> 1033.5385
$ws.Cells.Item(85, 12).Value = 1892   # L85: 2063.3333 -> 1892
$ws.Cells.Item(85, 13).Value = 214.4614999999999   # M85: 17 -> 214.4614999999999
$ws.Cells.Item(85, 14).Value = -4388   # N85: -4559.3333 -> -4388
$ws.Cells.Item(126, 8).Value = 4041.9023   # H126: 4074.1482 -> 4041.9023
$ws.Cells.Item(126, 9).Value = 3432.0178   # I126: 3449.8704 -> 3432.0178
$ws.Cells.Item(126, 10).Value = 5355.5   # J126: 5322.7036 -> 5355.5
$ws.Cells.Item(126, 11).Value = 10296.0534   # K126: 10349.6112 -> 10296.0534
$ws.Cells.Item(126, 12).Value = 16066.5   # L126: 15968.1108 -> 16066.5
$ws.Cells.Item(126, 13).Value = -7826.053400000001   # M126: -7879.611199999999 -> -7826.053400000001
$ws.Cells.Item(126, 14).Value = -21006.5   # N126: -20908.1108 -> -21006.5
$ws.Cells.Item(132, 8).Value = 2651.35   # H132: 2658.182 -> 2651.35
$ws.Cells.Item(132, 9).Value = 2604.2532   # I132: 2612.3206 -> 2604.2532
$ws.Cells.Item(132, 11).Value = 7812.7596   # K132: 7836.9618 -> 7812.7596
$ws.Cells.Item(132, 13).Value = -5282.7596   # M132: -5306.9618 -> -5282.7596
$ws.Cells.Item(139, 8).Value = 83707   # H139: 83707.5 -> 83707
$ws.Cells.Item(139, 10).Value = 83707   # J139: 83707.5 -> 83707
$ws.Cells.Item(139, 12).Value = 83707   # L139: 83707.5 -> 83707
$ws.Cells.Item(139, 14).Value = -93987   # N139: -93987.5 -> -93987

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(92, 8).Value = 50116.668   # H92: 50150 -> 50116.668
$ws.Cells.Item(92, 10).Value = 50116.668   # J92: 50150 -> 50116.668
$ws.Cells.Item(92, 12).Value = 50116.668   # L92: 50150 -> 50116.668
$ws.Cells.Item(92, 14).Value = -55108.668   # N92: -55142 -> -55108.668
$ws.Cells.Item(132, 8).Value = 1999.28   # H132: 1938.5385 -> 1999.28
$ws.Cells.Item(132, 9).Value = 800.375   # I132: 724.3 -> 800.375
$ws.Cells.Item(132, 11).Value = 2401.125   # K132: 2172.9 -> 2401.125
$ws.Cells.Item(132, 13).Value = 128.875   # M132: 357.1000000000004 -> 128.875
$ws.Cells.Item(138, 8).Value = 93119.8   # H138: 93140 -> 93119.8
$ws.Cells.Item(138, 10).Value = 91899.75   # J138: 91925 -> 91899.75
$ws.Cells.Item(138, 12).Value = 91899.75   # L138: 91925 -> 91899.75
$ws.Cells.Item(138, 14).Value = -102179.75   # N138: -102205 -> -102179.75
$ws.Cells.Item(141, 8).Value = 91809.664   # H141: 92357.5 -> 91809.664
$ws.Cells.Item(141, 10).Value = 91809.664   # J141: 92357.5 -> 91809.664
$ws.Cells.Item(141, 12).Value = 91809.664   # L141: 92357.5 -> 91809.664
$ws.Cells.Item(141, 14).Value = -102169.664   # N141: -102717.5 -> -102169.664
